$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the spicule type name in A2 from "Tylostyle" to "Oxea"
$ws.Range("A2").Value = "Oxea"

# Mimic natural Excel behavior: after editing A2 and pressing Enter,
# the active selection moves down to A3
$ws.Range("A3").Select()

$wb.Save()
